# Release dak-pnc 1.0.0 and pcmt-vax-prequal 0.2.0
# Update the pcmt-vaxprequal PreQualDB StructureDefinition metadata and
# remove the now-unused "Mapping: RIM Mapping" column from the Elements sheet.

$wb = $excel.ActiveWorkbook

# --- Metadata sheet: update Version, Status, Date, FHIR Version ---
$meta = $wb.Worksheets.Item("Metadata")

$meta.Range("B3").Value = "0.2.0"
$meta.Range("B6").Value = "active"
$meta.Range("B8").Value = "2025-09-16T20:42:07+00:00"
$meta.Range("B15").Value = "4.0.1"

# --- Elements sheet: delete column AK ("Mapping: RIM Mapping") ---
$elements = $wb.Worksheets.Item("Elements")
$elements.Columns.Item(37).Delete()
